$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the G column "profile" values for rows 18-41: "all" -> "LT"
for ($r = 18; $r -le 41; $r++) {
    $ws.Range("G$r").Value = "LT"
}

# 2. Clear the AutoFilter criteria (was filtering simtype == "red"),
#    which reveals the previously-hidden rows 2-17 again.
$ws.ShowAllData()

# 3. Update the active selection to match the new state.
$ws.Range("G18:G41").Select()
